$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Masthead text edits (rich-text runs inside A8 and C9).
#    Use Characters(start,len).Text so only the targeted substring is
#    replaced, matching the underlying shared-string run edits:
#      A8: "...Number  44"                -> "...Number  45"
#      C9: "...Week  10/31/2022  Through  11/6/2022"
#          -> "...Week  11/7/2022  Through  11/13/2022"
# ---------------------------------------------------------------------------

# A8 = "Volume 29   Number  44" -> replace the trailing "44" (chars 21-22)
$ws.Range("A8").Characters(21, 2).Text = "45"

# C9 = "Report Covering the Week  10/31/2022  Through  11/6/2022"
# Replace the second date first (chars 48-56) so the first replacement's
# length change doesn't shift the second date's position.
$ws.Range("C9").Characters(48, 9).Text = "11/13/2022"
$ws.Range("C9").Characters(27, 10).Text = "11/7/2022"

# ---------------------------------------------------------------------------
# 2) Weekly crime-stat table refresh (rows 14-29).
#    A handful of cells flip between a literal number and the sheet's
#    "insufficient data" placeholders (shared strings "0" / "***.*"), which
#    changes their stored type+style. Borrow the format from a donor cell
#    that already has the right combination, then (for number targets) set
#    the new numeric value.
# ---------------------------------------------------------------------------
$donorText0    = $ws.Range("D15")   # style 14, shared text "0"
$donorTextStar = $ws.Range("E15")   # style 14, shared text "***.*"
$donorNum15    = $ws.Range("F15")   # style 15, plain number
$donorNum16    = $ws.Range("H15")   # style 16, plain number

# --- cells whose storage type/style switches (text <-> number) ---
$donorText0.Copy($ws.Range("C15"))
$donorText0.Copy($ws.Range("C23"))
$donorText0.Copy($ws.Range("C26"))
$donorText0.Copy($ws.Range("D26"))
$donorTextStar.Copy($ws.Range("E26"))
$donorText0.Copy($ws.Range("D27"))
$donorTextStar.Copy($ws.Range("E27"))
$donorNum15.Copy($ws.Range("C28")); $ws.Range("C28").Value = 1
$donorNum15.Copy($ws.Range("D28")); $ws.Range("D28").Value = 1
$donorNum16.Copy($ws.Range("E28")); $ws.Range("E28").Value = 0
$donorNum15.Copy($ws.Range("F28")); $ws.Range("F28").Value = 1
$donorNum15.Copy($ws.Range("G28")); $ws.Range("G28").Value = 1
$donorNum16.Copy($ws.Range("H28")); $ws.Range("H28").Value = 0
$donorNum15.Copy($ws.Range("C29")); $ws.Range("C29").Value = 1
$donorNum15.Copy($ws.Range("D29")); $ws.Range("D29").Value = 1
$donorNum16.Copy($ws.Range("E29")); $ws.Range("E29").Value = 0
$donorNum15.Copy($ws.Range("F29")); $ws.Range("F29").Value = 1
$donorNum15.Copy($ws.Range("G29")); $ws.Range("G29").Value = 1
$donorNum16.Copy($ws.Range("H29")); $ws.Range("H29").Value = 0

# --- plain numeric value changes (style/type unchanged) ---
$values = @{
  "N14" = -78.378378378378
  "M15" = 45.454545454545
  "N15" = -54.929577464788
  "C16" = 2
  "D16" = 6
  "E16" = -66.666666666666
  "G16" = 14
  "H16" = -42.857142857142
  "I16" = 112
  "J16" = 132
  "K16" = -15.151515151515
  "L16" = -22.222222222222
  "M16" = -62.289562289562
  "N16" = -88.617886178861
  "C17" = 3
  "D17" = 9
  "E17" = -66.666666666666
  "G17" = 30
  "H17" = -26.666666666666
  "I17" = 337
  "J17" = 348
  "K17" = -3.160919540229
  "L17" = -8.423913043478
  "M17" = 9.415584415584
  "N17" = -49.550898203592
  "C18" = 1
  "D18" = 4
  "E18" = -75
  "G18" = 13
  "H18" = -23.076923076923
  "I18" = 111
  "J18" = 115
  "K18" = -3.478260869565
  "L18" = -19.565217391304
  "M18" = -67.543859649122
  "N18" = -88.497409326424
  "C19" = 8
  "D19" = 8
  "E19" = 0
  "F19" = 41
  "G19" = 32
  "H19" = 28.125
  "I19" = 390
  "J19" = 302
  "K19" = 29.139072847682
  "L19" = 10.481586402266
  "M19" = -31.338028169014
  "N19" = -88.760806916426
  "C20" = 2
  "D20" = 6
  "E20" = -66.666666666666
  "F20" = 14
  "G20" = 14
  "H20" = 0
  "I20" = 202
  "J20" = 156
  "K20" = 29.487179487179
  "L20" = 11.602209944751
  "M20" = -14.767932489451
  "N20" = -86.967741935483
  "C21" = 16
  "D21" = 33
  "E21" = -51.515151515151
  "F21" = 96
  "G21" = 105
  "H21" = -8.571428571428
  "I21" = 1192
  "J21" = 1085
  "K21" = 9.861751152073
  "L21" = -2.375102375102
  "M21" = -33.445002791736
  "N21" = -84.609425435765
  "M23" = -53.846153846153
  "C24" = 30
  "D24" = 29
  "E24" = 3.448275862068
  "F24" = 98
  "G24" = 90
  "H24" = 8.888888888888
  "I24" = 1154
  "J24" = 839
  "K24" = 37.544696066746
  "L24" = 28.794642857142
  "M24" = 22.896698615548
  "C25" = 16
  "D25" = 10
  "E25" = 60
  "F25" = 48
  "G25" = 43
  "H25" = 11.627906976744
  "I25" = 461
  "J25" = 430
  "K25" = 7.209302325581
  "L25" = 5.491990846681
  "M25" = -35.434173669467
  "G27" = 4
  "H27" = 75
  "I27" = 52
  "K27" = 20.930232558139
  "L27" = 136.363636363636
  "I28" = 31
  "J28" = 46
  "K28" = -32.608695652173
  "L28" = -24.390243902439
  "M28" = -39.215686274509
  "N28" = -78.620689655172
  "I29" = 25
  "J29" = 36
  "K29" = -30.555555555555
  "L29" = -26.470588235294
  "M29" = -40.476190476190
  "N29" = -80.620155038759
}
foreach ($key in $values.Keys) {
  $ws.Range($key).Value = $values[$key]
}
